$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New shared strings needed for this edit ---
# "T.J. Warren" and his bbref url already exist in the sheet (row 16) and will be
# reused in row 14. "Terrence Ross" and his url, "February 5, 1991" and
# "Washington" are brand-new values introduced by this edit (row 15).

# Row 14 currently holds Darius Bazley's data; it should become T.J. Warren's
# (his roster number/position/etc. move up one slot, and he gains a "No."
# value that wasn't present before).
$ws.Range("B14").Value = 21
$ws.Range("C14").Value = "T.J. Warren"
$ws.Range("D14").Value = "SF"
$ws.Range("E14").Value = "6-8"
$ws.Range("F14").Value = 220
$ws.Range("G14").Value = "September 5, 1993"
$ws.Range("H14").Value = "us"
$ws.Range("I14").Value = 7
$ws.Range("J14").Value = "NC State"
$ws.Range("K14").Value = "https://www.basketball-reference.com/players/w/warretj01.html"

# Row 15 currently holds Kevin Durant's data; it becomes a brand-new player,
# Terrence Ross, inserted into the roster.
$ws.Range("B15").Value = 8
$ws.Range("C15").Value = "Terrence Ross"
$ws.Range("D15").Value = "SG"
$ws.Range("E15").Value = "6-6"
$ws.Range("F15").Value = 206
$ws.Range("G15").Value = "February 5, 1991"
$ws.Range("H15").Value = "us"
$ws.Range("I15").Value = 10
$ws.Range("J15").Value = "Washington"
$ws.Range("K15").Value = "https://www.basketball-reference.com/players/r/rosste01.html"

# Row 16 currently holds T.J. Warren's data; it becomes Darius Bazley's data,
# which used to sit in row 14 (no roster "No." value, and no college listed).
$ws.Range("B16").ClearContents()
$ws.Range("C16").Value = "Darius Bazley"
$ws.Range("D16").Value = "PF"
$ws.Range("E16").Value = "6-8"
$ws.Range("F16").Value = 208
$ws.Range("G16").Value = "June 12, 2000"
$ws.Range("H16").Value = "us"
$ws.Range("I16").Value = 3
$ws.Range("J16").ClearContents()
$ws.Range("K16").Value = "https://www.basketball-reference.com/players/b/bazleda01.html"

# Row 17 is a brand-new row holding Kevin Durant's data, which used to sit in
# row 15. Copy formatting down from row 16 first so the new row matches the
# rest of the table (bold/bordered "No." cell, hyperlink style on the url).
$ws.Range("A16:K16").Copy()
$ws.Range("A17:K17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A17").Value = 15
$ws.Range("C17").Value = "Kevin Durant"
$ws.Range("D17").Value = "SF"
$ws.Range("E17").Value = "6-10"
$ws.Range("F17").Value = 240
$ws.Range("G17").Value = "September 29, 1988"
$ws.Range("H17").Value = "us"
$ws.Range("I17").Value = 14
$ws.Range("J17").Value = "Texas"
$ws.Range("K17").Value = "https://www.basketball-reference.com/players/d/duranke01.html"

# Register the actual hyperlink relationships for the two urls that changed
# cell (row 14 now points at T.J. Warren's page, row 17 is brand new).
$ws.Hyperlinks.Add($ws.Range("K14"), "https://www.basketball-reference.com/players/w/warretj01.html") | Out-Null
$ws.Hyperlinks.Add($ws.Range("K15"), "https://www.basketball-reference.com/players/r/rosste01.html") | Out-Null
$ws.Hyperlinks.Add($ws.Range("K17"), "https://www.basketball-reference.com/players/d/duranke01.html") | Out-Null

# Adding a hyperlink re-styles the cell with a fresh (unnamed) style variant;
# reset it back to the workbook's existing "Hyperlink" cell style so K14:K17
# all share the same style index the rest of the K column already uses.
$ws.Range("K14").Style = "Hyperlink"
$ws.Range("K15").Style = "Hyperlink"
$ws.Range("K16").Style = "Hyperlink"
$ws.Range("K17").Style = "Hyperlink"
